$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.759.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -7.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.688.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -7.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.677.88"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.624"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -10.14%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.703"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -12.79%  "
$ws.Range("E11").Value = "  -13.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "51.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -8.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000292"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -13.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -11.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.259.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.689.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -10.37%  "
$ws.Range("E18").Value = "  -3.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -10.00%  "
$ws.Range("E20").Value = "  -10.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.699.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "407.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -11.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -9.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -10.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -11.04%  "
$ws.Range("E27").Value = "  -5.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -12.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -11.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.19%  "
$ws.Range("E32").Value = "  -5.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -11.52%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.116"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -10.21%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "64.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.67%  "
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "43.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -12.31%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "602.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0888"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -14.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.397"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("E42").Value = "  -10.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.98"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -12.67%  "
$ws.Range("E44").Value = "  -11.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0436"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -10.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -13.35%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.778.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.68%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.20%  "
$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.55%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.134"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -11.34%  "
